$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2) | Out-Null
}

Replace-Text "2024-12-18 Wednesday" "2024-12-19 Thursday"

Replace-Text "46×23=" "91×26="
Replace-Text "93×18=" "14×56="
Replace-Text "81×39=" "80×32="
Replace-Text "91×54=" "93×35="
Replace-Text "21×68=" "67×84="

Replace-Text "45×23=" "65×95="
Replace-Text "38×77=" "56×27="
Replace-Text "14×55=" "24×62="
Replace-Text "39×65=" "42×23="
Replace-Text "33×48=" "65×66="

Replace-Text "75×94=" "89×86="
Replace-Text "22×46=" "75×93="
Replace-Text "97×84=" "18×58="
Replace-Text "16×64=" "54×28="
Replace-Text "51×44=" "95×58="

Replace-Text "34×86=" "15×11="
Replace-Text "94×52=" "40×12="
Replace-Text "13×58=" "41×64="
Replace-Text "72×30=" "57×22="
Replace-Text "81×89=" "81×94="

Replace-Text "75×21=" "47×66="
Replace-Text "78×96=" "89×49="
Replace-Text "27×48=" "80×38="
Replace-Text "14×21=" "68×60="
Replace-Text "42×18=" "94×38="
